# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style used by the other header cells (B1:H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy the formatting of the existing header cell (H1) so the
# new header cells share the same style (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for rows 2-40 (row index -> I, J)
$iVals = @(9,9,9,8,8,9,8,9,10,9,8,8,8,7,7,8,7,7,7,5,9,10,9,8,8,6,1,6,6,8,8,6,6,5,8,4,4,4,5)
$jVals = @(9,9,9,9,8,9,9,9,10,9,9,8,8,8,8,8,7,7,7,6,9,10,9,8,8,7,2,7,7,8,8,6,6,5,8,4,4,4,5)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
